$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.838.10'
$ws.Range("E2").Value = '  -2.01%  '

$ws.Range("D3").Value = '3.316.63'
$ws.Range("E3").Value = '  +1.37%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = "'574.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.66%  '

$ws.Range("D6").Value = "'182.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.20%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").Value = "'0.601"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").Value = '3.319.89'
$ws.Range("E9").Value = '  +1.49%  '

$ws.Range("D10").Value = "'0.128"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.04%  '

$ws.Range("D11").Value = "'6.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.02%  '

$ws.Range("D12").Value = "'0.404"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.91%  '

$ws.Range("D13").Value = '3.900.60'
$ws.Range("E13").Value = '  +1.54%  '

$ws.Range("E14").Value = '  -0.89%  '

$ws.Range("D15").Value = "'27.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.99%  '

$ws.Range("D16").Value = '67.021.67'
$ws.Range("E16").Value = '  -1.67%  '

$ws.Range("E17").Value = '  -0.80%  '

$ws.Range("D18").Value = '3.342.44'
$ws.Range("E18").Value = '  +2.14%  '

$ws.Range("D19").Value = "'438.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.47%  '

$ws.Range("D20").Value = "'13.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.69%  '

$ws.Range("D21").Value = "'5.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.36%  '

$ws.Range("D22").Value = "'7.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.03%  '

$ws.Range("D23").Value = "'73.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.68%  '

$ws.Range("D24").Value = "'0.997"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.50%  '

$ws.Range("D25").Value = '3.475.93'
$ws.Range("E25").Value = '  +1.83%  '

$ws.Range("D26").Value = "'0.510"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.34%  '

$ws.Range("D27").Value = "'0.0000118"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.78%  '

$ws.Range("E28").Value = '  +3.65%  '

$ws.Range("D29").Value = "'8.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.66%  '

$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.54%  '

$ws.Range("D31").Value = "'1.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.61%  '

$ws.Range("D32").Value = "'22.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.99%  '

$ws.Range("E33").Value = '  -0.08%  '

$ws.Range("D34").Value = "'5.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.98%  '

$ws.Range("D35").Value = "'6.76"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.09%  '

$ws.Range("E36").Value = '  -2.30%  '

$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = "'1.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.74%  '

$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").Value = "'161.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.82%  '

$ws.Range("D39").Value = "'27.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.26%  '

$ws.Range("D40").Value = "'1.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.65%  '

$ws.Range("D41").Value = '2.815.01'
$ws.Range("E41").Value = '  +6.91%  '

$ws.Range("D42").Value = "'0.786"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.56%  '

$ws.Range("D43").Value = "'4.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.19%  '

$ws.Range("D44").Value = "'6.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.64%  '

$ws.Range("D45").Value = "'40.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.78%  '

$ws.Range("D46").Value = "'0.0670"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.79%  '

$ws.Range("D47").Value = "'24.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.04%  '

$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").Value = "'2.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.68%  '

$ws.Range("B49").Value = 'Bittensor'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D49").Value = "'320.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.34%  '

$ws.Range("D50").Value = "'0.0271"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.87%  '

$ws.Range("D51").Value = "'0.977"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.04%  '
